$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 58 ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(58, 1).Value = "Ruilen van product"
$logs.Cells.Item(58, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(58, 3).Value = "Kan ik dit product ruilen voor een andere maat?"
$logs.Cells.Item(58, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(58, 6).Value = "2025-06-22 22:13:34"
$logs.Cells.Item(58, 7).Value = "Nee"

# Extend the conditional formatting ranges to cover the new row 58
$catFcs = $logs.Range("D2:D57").FormatConditions
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($logs.Range("D2:D58"))
}

$answeredFcs = $logs.Range("G2:G57").FormatConditions
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($logs.Range("G2:G58"))
}

# --- Dashboard sheet: bump "Retour / Terugbetaling" count from 7 to 8 ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 8
